# Update "想去人数" (want-to-go count) values in column F for the
# 展览 (Exhibitions) sheet and the 全部类型 (All types) sheet, reflecting
# the regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 7978
$ws1.Range("F10").Value = 468
$ws1.Range("F13").Value = 452
$ws1.Range("F14").Value = 69
$ws1.Range("F15").Value = 76
$ws1.Range("F17").Value = 5873
$ws1.Range("F18").Value = 184
$ws1.Range("F19").Value = 266
$ws1.Range("F20").Value = 1843
$ws1.Range("F21").Value = 8
$ws1.Range("F22").Value = 14
$ws1.Range("F23").Value = 239
$ws1.Range("F24").Value = 396

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 7979
$ws4.Range("F10").Value = 468
$ws4.Range("F13").Value = 452
$ws4.Range("F14").Value = 69
$ws4.Range("F15").Value = 76
$ws4.Range("F18").Value = 5873
$ws4.Range("F20").Value = 184
$ws4.Range("F21").Value = 266
$ws4.Range("F22").Value = 1843
$ws4.Range("F23").Value = 8
$ws4.Range("F24").Value = 14
$ws4.Range("F25").Value = 239
$ws4.Range("F26").Value = 396

$wb.Save()
